$d = $word.ActiveDocument

# Locate the title paragraph ("            CS 1632 Software Quality Assurance")
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "CS\s*1632") {
        # Center the title paragraph.
        $p.Format.Alignment = 1

        # Remove the stray leading run of spaces that precedes "CS 1632...".
        $paraText = $p.Range.Text
        $trimmed = $paraText.TrimStart(" ")
        $leadingCount = $paraText.Length - $trimmed.Length

        if ($leadingCount -gt 0) {
            $pStart = $p.Range.Start
            $delRng = $d.Range($pStart, $pStart + $leadingCount)
            $delRng.Delete()
        }

        break
    }
}
